$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-32: update date serial value from 45221 to 45224
$ws.Range("C2:C32").Value = 45224
